$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 28 and 29 (pushing the U1Rx/U1Tx rows down from 29/30 to 31/32)
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(29).Insert()

# Fill in the B/F values for the two new rows for the front IR sensors
$ws.Range("B28").Value = 6
$ws.Range("F28").Value = "5V"
$ws.Range("B29").Value = 7
$ws.Range("F29").Value = "5V"

# Update the G-column labels (camera -> IR sensor naming). Written in this
# order so the new shared-string table entries land in the same order as
# the source workbook.
$ws.Range("G29").Value = "Front Right IR"
$ws.Range("G28").Value = "Front Left IR"
$ws.Range("G27").Value = "back IR"
$ws.Range("G26").Value = "right IR"
$ws.Range("G25").Value = "left IR"

# Update the selection / active cell (was G31, now G27) and clear the scrolled
# top-left cell so the view resets to the top of the sheet
$ws.Range("A1").Select()
$ws.Range("G27").Select()
